$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 45455000
$ws.Range("I6").Value = 71429060
$ws.Range("K6").Value = 214287180
$ws.Range("M6").Value = -214287068

$ws.Range("H31").Value = 300
$ws.Range("I31").Value = 300
$ws.Range("K31").Value = 900
$ws.Range("M31").Value = -670

$ws.Range("H33").Value = 1837.3684
$ws.Range("I33").Value = 252.91667
$ws.Range("J33").Value = 4553.5713
$ws.Range("K33").Value = 252.91667
$ws.Range("L33").Value = 4553.5713
$ws.Range("M33").Value = -23.91667000000001
$ws.Range("N33").Value = -5011.5713

$ws.Range("H39").Value = 193.66667
$ws.Range("J39").Value = 1000
$ws.Range("L39").Value = 3000
$ws.Range("N39").Value = -3592

$ws.Range("H43").Value = 19826
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 19826
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 19826
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -19964

$ws.Range("H98").Value = 1095.2
$ws.Range("J98").Value = 1199.3334
$ws.Range("L98").Value = 1199.3334
$ws.Range("N98").Value = -4195.3334

$ws.Range("H100").Value = 5057.4287
$ws.Range("I100").Value = 2222
$ws.Range("K100").Value = 2222
$ws.Range("M100").Value = -1681

$ws.Range("H122").Value = 1095.2
$ws.Range("J122").Value = 1199.3334
$ws.Range("L122").Value = 3598.0002
$ws.Range("N122").Value = -8498.0002

$ws.Range("H129").Value = 2540.92
$ws.Range("J129").Value = 2968.9473
$ws.Range("L129").Value = 8906.841899999999
$ws.Range("N129").Value = -18906.8419

$ws.Range("H137").Value = 2222
$ws.Range("I137").Value = 1579.7
$ws.Range("K137").Value = 4739.1
$ws.Range("M137").Value = -2189.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2675
$ws.Range("I32").Value = 2241.3157
$ws.Range("K32").Value = 2241.3157
$ws.Range("M32").Value = -1954.3157

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 569.2963
$ws.Range("J80").Value = 650.82355
$ws.Range("L80").Value = 650.82355
$ws.Range("N80").Value = -2646.82355

$ws.Range("H83").Value = 569.2963
$ws.Range("J83").Value = 650.82355
$ws.Range("L83").Value = 3254.11775
$ws.Range("N83").Value = -13238.11775

$ws.Range("H86").Value = 3131.5625
$ws.Range("I86").Value = 2300.6
$ws.Range("J86").Value = 4516.5
$ws.Range("K86").Value = 2300.6
$ws.Range("L86").Value = 4516.5
$ws.Range("M86").Value = -1177.6
$ws.Range("N86").Value = -6762.5

$ws.Range("H89").Value = 3131.5625
$ws.Range("I89").Value = 2300.6
$ws.Range("J89").Value = 4516.5
$ws.Range("K89").Value = 11503
$ws.Range("L89").Value = 22582.5
$ws.Range("M89").Value = -5887
$ws.Range("N89").Value = -33814.5

$ws.Range("H105").Value = 4037.5
$ws.Range("I105").Value = 3566
$ws.Range("K105").Value = 3566
$ws.Range("M105").Value = -1819

$ws.Range("H107").Value = 7291.125
$ws.Range("I107").Value = 5332.375
$ws.Range("J107").Value = 9249.875
$ws.Range("K107").Value = 5332.375
$ws.Range("L107").Value = 9249.875
$ws.Range("M107").Value = -3412.375
$ws.Range("N107").Value = -13089.875

$ws.Range("H134").Value = 2475.606
$ws.Range("I134").Value = 2370.762
$ws.Range("J134").Value = 2659.0833
$ws.Range("K134").Value = 7112.286
$ws.Range("L134").Value = 7977.249899999999
$ws.Range("M134").Value = -4577.286
$ws.Range("N134").Value = -13047.2499

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2659.8
$ws.Range("J58").Value = 1680
$ws.Range("L58").Value = 1680
$ws.Range("N58").Value = -2086

$ws.Range("H122").Value = 1844.4166
$ws.Range("I122").Value = 1648.909
$ws.Range("J122").Value = 3995
$ws.Range("K122").Value = 4946.727000000001
$ws.Range("L122").Value = 11985
$ws.Range("M122").Value = -2496.727000000001
$ws.Range("N122").Value = -16885

$ws.Range("H136").Value = 2659.8
$ws.Range("J136").Value = 1680
$ws.Range("L136").Value = 5040
$ws.Range("N136").Value = -10140

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 922.75
$ws.Range("J38").Value = 1372.5
$ws.Range("L38").Value = 4117.5
$ws.Range("N38").Value = -4811.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5937.25
$ws.Range("I80").Value = 3099.5
$ws.Range("K80").Value = 3099.5
$ws.Range("M80").Value = -2101.5

$ws.Range("H83").Value = 5937.25
$ws.Range("I83").Value = 3099.5
$ws.Range("K83").Value = 15497.5
$ws.Range("M83").Value = -10505.5

$ws.Range("H132").Value = 2852.611
$ws.Range("I132").Value = 2352.3845
$ws.Range("J132").Value = 4153.2
$ws.Range("K132").Value = 7057.1535
$ws.Range("L132").Value = 12459.6
$ws.Range("M132").Value = -4527.1535
$ws.Range("N132").Value = -17519.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1645
$ws.Range("I7").Value = 1526.8334
$ws.Range("J7").Value = 1999.5
$ws.Range("K7").Value = 1526.8334
$ws.Range("L7").Value = 1999.5
$ws.Range("M7").Value = -1414.8334
$ws.Range("N7").Value = -2223.5

$ws.Range("H22").Value = 2177.8108
$ws.Range("I22").Value = 1656.4546
$ws.Range("K22").Value = 1656.4546
$ws.Range("M22").Value = -1361.4546

$ws.Range("H27").Value = 2177.8108
$ws.Range("I27").Value = 1656.4546
$ws.Range("K27").Value = 1656.4546
$ws.Range("M27").Value = -1549.4546

$ws.Range("H46").Value = 3316
$ws.Range("I46").Value = 2229
$ws.Range("J46").Value = 3436.7778
$ws.Range("K46").Value = 2229
$ws.Range("L46").Value = 3436.7778
$ws.Range("M46").Value = -2041
$ws.Range("N46").Value = -3812.7778

$ws.Range("H68").Value = 7126.533
$ws.Range("J68").Value = 8299.9
$ws.Range("L68").Value = 8299.9
$ws.Range("N68").Value = -9797.9

$ws.Range("H71").Value = 7126.533
$ws.Range("J71").Value = 8299.9
$ws.Range("L71").Value = 41499.5
$ws.Range("N71").Value = -48987.5

$ws.Range("H93").Value = 3115.3635
$ws.Range("I93").Value = 633.9048
$ws.Range("K93").Value = 633.9048
$ws.Range("M93").Value = 614.0952

$ws.Range("H126").Value = 1645
$ws.Range("I126").Value = 1526.8334
$ws.Range("J126").Value = 1999.5
$ws.Range("K126").Value = 4580.5002
$ws.Range("L126").Value = 5998.5
$ws.Range("M126").Value = -2110.5002
$ws.Range("N126").Value = -10938.5

$ws.Range("H132").Value = 2426.9167
$ws.Range("I132").Value = 2374.3333
$ws.Range("K132").Value = 7122.999899999999
$ws.Range("M132").Value = -4592.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 948.2222
$ws.Range("I81").Value = 948.2222
$ws.Range("K81").Value = 1896.4444
$ws.Range("M81").Value = -835.4444000000001

$ws.Range("H84").Value = 948.2222
$ws.Range("I84").Value = 948.2222
$ws.Range("K84").Value = 9482.222
$ws.Range("M84").Value = -4178.222

$ws.Range("H122").Value = 5893.5
$ws.Range("I122").Value = 3412
$ws.Range("K122").Value = 10236
$ws.Range("M122").Value = -7786

$ws.Range("H132").Value = 5475.625
